$d = $word.ActiveDocument
$d.Bookmarks.ShowHidden = $true

# ---------------------------------------------------------------
# 1) "TBD" (paragraph 4, under "Confirmation of Deliverables")
#    -> "3 index cards for each member", then add three new
#    Heading2 paragraphs after it.
# ---------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Text = "3 index cards for each member"

$r = $p4.Range
$r.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Text = "Task Board Review"

$r = $p5.Range
$r.InsertParagraphAfter()
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Text = "Updated Burn Down"

$r = $p6.Range
$r.InsertParagraphAfter()
$p7 = $d.Paragraphs.Item(7)
$p7.Range.Text = "Iteration 2 Demonstration Summary"

# ---------------------------------------------------------------
# 2) "Review last weeks minutes" -> "Review last week" + right
#    single quote + "s minutes"  (paragraph shifted down by 3,
#    now item 9)
# ---------------------------------------------------------------
$p9 = $d.Paragraphs.Item(9)
$start = $p9.Range.Start
$p9.Range.Text = "Review last week"
$afterFirst = $start + ("Review last week".Length)
$ins1 = $d.Range($afterFirst, $afterFirst)
$ins1.InsertAfter([char]0x2019)
$afterSecond = $afterFirst + 1
$ins2 = $d.Range($afterSecond, $afterSecond)
$ins2.InsertAfter("s minutes")

# ---------------------------------------------------------------
# 3) Remove the "Progress Report" paragraph entirely (it carried
#    an explicit black-color run/paragraph formatting that is not
#    present anywhere in the new text).  This paragraph is now
#    item 13.
# ---------------------------------------------------------------
$p13 = $d.Paragraphs.Item(13)
$p13.Range.Delete()

# ---------------------------------------------------------------
# 4) "Task Assignments" (now item 14, after the deletion above)
#    gains a trailing space run: "Task Assignments "
# ---------------------------------------------------------------
$p14 = $d.Paragraphs.Item(14)
$start14 = $p14.Range.Start
$p14.Range.Text = "Task Assignments"
$afterTask = $start14 + ("Task Assignments".Length)
$ins3 = $d.Range($afterTask, $afterTask)
$ins3.InsertAfter(" ")

# ---------------------------------------------------------------
# 5) The old "TBD" paragraph (now item 17, following "New
#    Business" / "Iteration 1 Assessment" which are unchanged)
#    becomes "Graphical Assessment of Apps", and the _GoBack
#    bookmark moves from the end of the paragraph to the start.
# ---------------------------------------------------------------
$p17 = $d.Paragraphs.Item(17)
$start17 = $p17.Range.Start
$p17.Range.Text = "Graphical Assessment of Apps"

$d.Bookmarks.Item("_GoBack").Delete()
$bmRange = $d.Range($start17, $start17)
$d.Bookmarks.Add("_GoBack", $bmRange)
